$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns remain text so values like "1.70" or "64.30"
# keep their exact formatting instead of being auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '44.389.33'
$ws.Range("E2").Value = '  +0.35%  '
$ws.Range("D3").Value = '2.217.49'
$ws.Range("E3").Value = '  -0.87%  '
$ws.Range("E4").Value = '  +0.46%  '
$ws.Range("D5").Value = '301.47'
$ws.Range("E5").Value = '  -1.75%  '
$ws.Range("D6").Value = '89.11'
$ws.Range("E6").Value = '  -4.90%  '
$ws.Range("D7").Value = '0.558'
$ws.Range("E7").Value = '  -2.14%  '
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("D9").Value = '0.496'
$ws.Range("E9").Value = '  -4.87%  '
$ws.Range("D10").Value = '33.59'
$ws.Range("E10").Value = '  -2.08%  '
$ws.Range("D11").Value = '0.0778'
$ws.Range("E11").Value = '  -3.45%  '
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").Value = '0.104'
$ws.Range("E12").Value = '  -0.47%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '6.92'
$ws.Range("E13").Value = '  -2.89%  '
$ws.Range("D14").Value = '2.566.37'
$ws.Range("E14").Value = '  -0.47%  '
$ws.Range("D15").Value = '2.226.55'
$ws.Range("E15").Value = '  -4.18%  '
$ws.Range("D16").Value = '0.803'
$ws.Range("E16").Value = '  -2.49%  '
$ws.Range("D17").Value = '13.14'
$ws.Range("E17").Value = '  -1.92%  '
$ws.Range("D18").Value = '44.376.91'
$ws.Range("E18").Value = '  +0.92%  '
$ws.Range("D19").Value = '0.0₃0909'
$ws.Range("E19").Value = '  -5.40%  '
$ws.Range("D20").Value = '6.01'
$ws.Range("E20").Value = '  -5.09%  '
$ws.Range("D21").Value = '11.28'
$ws.Range("E21").Value = '  -6.24%  '
$ws.Range("D22").Value = '64.30'
$ws.Range("E22").Value = '  -1.79%  '
$ws.Range("D23").Value = '232.79'
$ws.Range("E23").Value = '  -1.23%  '
$ws.Range("D24").Value = '2.86'
$ws.Range("E24").Value = '  -2.56%  '
$ws.Range("E25").Value = '  -0.01%  '
$ws.Range("D26").Value = '1.91'
$ws.Range("E26").Value = '  -4.48%  '
$ws.Range("D27").Value = '2.26'
$ws.Range("E27").Value = '  +2.66%  '
$ws.Range("D28").Value = '9.45'
$ws.Range("E28").Value = '  -3.85%  '
$ws.Range("D29").Value = '36.09'
$ws.Range("E29").Value = '  -8.17%  '
$ws.Range("D30").Value = '19.46'
$ws.Range("E30").Value = '  -2.50%  '
$ws.Range("D31").Value = '5.55'
$ws.Range("E31").Value = '  -4.47%  '
$ws.Range("D32").Value = '145.42'
$ws.Range("E32").Value = '  -5.12%  '
$ws.Range("E33").Value = '  +0.44%  '
$ws.Range("D34").Value = '0.0754'
$ws.Range("E34").Value = '  -4.63%  '
$ws.Range("D35").Value = '2.98'
$ws.Range("E35").Value = '  -3.59%  '
$ws.Range("D36").Value = '0.106'
$ws.Range("E36").Value = '  -1.44%  '
$ws.Range("E37").Value = '  -3.45%  '
$ws.Range("D38").Value = '1.77'
$ws.Range("E38").Value = '  +2.21%  '
$ws.Range("D39").Value = '14.41'
$ws.Range("E39").Value = '  +1.94%  '
$ws.Range("D40").Value = '3.22'
$ws.Range("E40").Value = '  -6.96%  '
$ws.Range("D41").Value = '3.62'
$ws.Range("E41").Value = '  -4.40%  '
$ws.Range("D42").Value = '0.0286'
$ws.Range("E42").Value = '  -3.69%  '
$ws.Range("E43").Value = '  +0.20%  '
$ws.Range("D44").Value = '1.774.74'
$ws.Range("E44").Value = '  +2.96%  '
$ws.Range("D45").Value = '1.70'
$ws.Range("E45").Value = '  +6.64%  '
$ws.Range("D46").Value = '78.82'
$ws.Range("E46").Value = '  -3.96%  '
$ws.Range("D47").Value = '0.181'
$ws.Range("E47").Value = '  -5.02%  '
$ws.Range("D48").Value = '95.16'
$ws.Range("E48").Value = '  -3.79%  '
$ws.Range("D49").Value = '4.72'
$ws.Range("E49").Value = '  -3.88%  '
$ws.Range("D50").Value = '66.43'
$ws.Range("E50").Value = '  -0.62%  '
$ws.Range("D51").Value = '52.14'
$ws.Range("E51").Value = '  -4.16%  '
